# Nathan's hours wk 7
# Fill in the week-7 booked-hours rows (23-28) on Sheet1 and tag the blank
# trailing row (29) with the week number, matching the source timesheet
# entries added in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-Activity($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.Value = $text
    # Match the font used by the other populated "Activity" cells in this
    # column (Verdana 10, unprotected) so the new rows look like the rest
    # of the table.
    $r.Font.Name = "Verdana"
    $r.Font.Size = 10
}

# Row 23: 04/09/2023 09:30 - 10:00
$ws.Range("A23").Value = 7
$ws.Range("B23").Value = 43711
$ws.Range("C23").Value = 0.39583333333333331
$ws.Range("D23").Value = 43711
$ws.Range("E23").Value = 0.41666666666666669
Set-Activity "F23" "Team meeting (study break)"

# Row 24: 07/09/2023 17:00 - 21:00
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = 43714
$ws.Range("C24").Value = 0.70833333333333337
$ws.Range("D24").Value = 43714
$ws.Range("E24").Value = 0.875
Set-Activity "F24" "Set up Auth0 app & integrated API"

# Row 25: 08/09/2023 20:00 - 21:00
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = 43715
$ws.Range("C25").Value = 0.83333333333333337
$ws.Range("D25").Value = 43715
$ws.Range("E25").Value = 0.875
Set-Activity "F25" "Set up basic login button template, created routing for profile access"

# Row 26: 11/09/2023 10:30 - 11:00
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = 43718
$ws.Range("C26").Value = 0.4375
$ws.Range("D26").Value = 43718
$ws.Range("E26").Value = 0.45833333333333331
Set-Activity "F26" "Formal team meeting"

# Row 27: 11/09/2023 11:30 - 12:00
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = 43718
$ws.Range("C27").Value = 0.47916666666666669
$ws.Range("D27").Value = 43718
$ws.Range("E27").Value = 0.5
Set-Activity "F27" "Client briefing/meeting"

# Row 28: 15/09/2023 20:00 - 22:00
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = 43722
$ws.Range("C28").Value = 0.83333333333333337
$ws.Range("D28").Value = 43722
$ws.Range("E28").Value = 0.91666666666666663
Set-Activity "F28" "Testing Auth0 authorisation methods"

# Row 29 stays blank except for the week-number tag carried into the next
# (still empty) entry row.
$ws.Range("A29").Value = 7

# Leave the view/selection near the newly-entered data, mirroring the
# author's last on-screen position (scrolled down, zoomed to 137%, E25
# selected).
$excel.ActiveWindow.Zoom = 137
$ws.Range("E25").Select() | Out-Null
